# RPA-93: Lage initieringsfil og rutine
# Insert a new "Kommune_Nr" column right after "Sak_Nr" (i.e. before the
# existing "FNR" column), fill in its sample value, and populate two
# cells on the sample row that were previously left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting a whole column at C shifts every existing column from C
# onward one position to the right (C->D, D->E, ... AB->AC), which is
# exactly what the diff shows happening to the rest of the row.
$ws.Columns("C:C").Insert()

# New header + sample value for the inserted "Kommune_Nr" column.
$ws.Range("C1").Value = "Kommune_Nr"
$ws.Range("C2").Value = 11111

# "_2_barn1118" (now column L after the shift) gets an explicit 0 value.
$ws.Range("L2").Value = 0

# "_4_Status" (now column AB after the shift) gets the "OK" text. It
# already carried the wrap-text style (s="1") before the edit, so
# re-apply WrapText to keep that style on the cell after writing the
# value.
$ws.Range("AB2").Value = "OK"
$ws.Range("AB2").WrapText = $true
